$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.069.58'
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("D3").Value = '3.406.85'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  -0.03%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '572.52'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.14%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '163.14'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.43%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.408.76'
$ws.Range("E8").Value = '  -1.08%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.551'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -4.65%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.31'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.58%  '
$ws.Range("E11").Value = '  -1.69%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.423'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -4.08%  '
$ws.Range("D13").Value = '3.987.67'
$ws.Range("E13").Value = '  -1.41%  '
$ws.Range("E14").Value = '  +0.69%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '27.04'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.31%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.0000173'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.64%  '
$ws.Range("D17").Value = '64.069.16'
$ws.Range("E17").Value = '  -1.37%  '
$ws.Range("D18").Value = '3.390.94'
$ws.Range("E18").Value = '  -1.94%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.14'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.36%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.63'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.67%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '377.97'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.29%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '7.79'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.11%  '
$ws.Range("E23").Value = '  +0.20%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '70.20'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.77%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.512'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -4.99%  '
$ws.Range("E26").Value = '  -4.26%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.49'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -3.79%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.179'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.23%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.25%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '6.13'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.13%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.40'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -3.74%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.01'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.22%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '22.88'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -1.49%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '7.09'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.47%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.49'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -4.89%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '160.07'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.46%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.860'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +9.31%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.81'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -3.30%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0723'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -3.43%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '42.86'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '2.763.89'
$ws.Range("E41").Value = '  -5.01%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '25.68'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.31%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '26.30'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.95%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '6.44'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.88%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '4.39'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -3.06%  '
$ws.Range("E46").Value = '  -1.59%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.44'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +4.24%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '329.14'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.81%  '
$ws.Range("E49").Value = '  -3.76%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.102'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.82%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '6.28'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.96%  '